$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 60630.11
$ws.Range("J17").Value = 62300.39
$ws.Range("L17").Value = 186901.17
$ws.Range("N17").Value = -187237.17
$ws.Range("H21").Value = 22901.6
$ws.Range("I21").Value = 35009.5
$ws.Range("J21").Value = 19874.625
$ws.Range("K21").Value = 35009.5
$ws.Range("L21").Value = 19874.625
$ws.Range("M21").Value = -34541.5
$ws.Range("N21").Value = -20810.625
$ws.Range("H23").Value = 22901.6
$ws.Range("I23").Value = 35009.5
$ws.Range("J23").Value = 19874.625
$ws.Range("K23").Value = 35009.5
$ws.Range("L23").Value = 19874.625
$ws.Range("M23").Value = -34775.5
$ws.Range("N23").Value = -20342.625
$ws.Range("H40").Value = 1961.9656
$ws.Range("I40").Value = 1936.2273
$ws.Range("J40").Value = 2042.8572
$ws.Range("K40").Value = 1936.2273
$ws.Range("L40").Value = 2042.8572
$ws.Range("M40").Value = -1761.2273
$ws.Range("N40").Value = -2392.8572
$ws.Range("H76").Value = 3174
$ws.Range("I76").Value = 3109.75
$ws.Range("J76").Value = 3302.5
$ws.Range("K76").Value = 3109.75
$ws.Range("L76").Value = 3302.5
$ws.Range("M76").Value = -2794.75
$ws.Range("N76").Value = -3932.5
$ws.Range("H79").Value = 3174
$ws.Range("I79").Value = 3109.75
$ws.Range("J79").Value = 3302.5
$ws.Range("K79").Value = 3109.75
$ws.Range("L79").Value = 3302.5
$ws.Range("M79").Value = -2017.75
$ws.Range("N79").Value = -5486.5
$ws.Range("H98").Value = 2804
$ws.Range("I98").Value = 2516.16
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 2516.16
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = -1018.16
$ws.Range("N98").Value = -12996
$ws.Range("H122").Value = 2804
$ws.Range("I122").Value = 2516.16
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 7548.48
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -5098.48
$ws.Range("N122").Value = -34900
$ws.Range("H129").Value = 5953570.5
$ws.Range("J129").Value = 1098.775
$ws.Range("L129").Value = 3296.325
$ws.Range("N129").Value = -13296.325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3261.1765
$ws.Range("I61").Value = 2725.6
$ws.Range("K61").Value = 2725.6
$ws.Range("M61").Value = -2513.6
$ws.Range("H119").Value = 1285611
$ws.Range("J119").Value = 1285611
$ws.Range("L119").Value = 1285611
$ws.Range("N119").Value = -1295287
$ws.Range("H132").Value = 2007.7455
$ws.Range("I132").Value = 1537
$ws.Range("K132").Value = 4611
$ws.Range("M132").Value = -2081
$ws.Range("H136").Value = 3261.1765
$ws.Range("I136").Value = 2725.6
$ws.Range("K136").Value = 8176.799999999999
$ws.Range("M136").Value = -5626.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1724.2142
$ws.Range("J105").Value = 3000
$ws.Range("L105").Value = 3000
$ws.Range("N105").Value = -6494
$ws.Range("H134").Value = 6083.385
$ws.Range("I134").Value = 5889.304
$ws.Range("J134").Value = 7571.3335
$ws.Range("K134").Value = 17667.912
$ws.Range("L134").Value = 22714.0005
$ws.Range("M134").Value = -15132.912
$ws.Range("N134").Value = -27784.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1410.75
$ws.Range("I16").Value = 774.7692
$ws.Range("J16").Value = 4166.6665
$ws.Range("K16").Value = 774.7692
$ws.Range("L16").Value = 4166.6665
$ws.Range("M16").Value = -487.7692
$ws.Range("N16").Value = -4740.6665
$ws.Range("H58").Value = 17862116
$ws.Range("I58").Value = 3666.7856
$ws.Range("K58").Value = 3666.7856
$ws.Range("M58").Value = -3463.7856
$ws.Range("H99").Value = 5719.2856
$ws.Range("I99").Value = 2740.3333
$ws.Range("J99").Value = 7953.5
$ws.Range("K99").Value = 2740.3333
$ws.Range("L99").Value = 7953.5
$ws.Range("M99").Value = -1242.3333
$ws.Range("N99").Value = -10949.5
$ws.Range("H113").Value = 1410.75
$ws.Range("I113").Value = 774.7692
$ws.Range("J113").Value = 4166.6665
$ws.Range("K113").Value = 774.7692
$ws.Range("L113").Value = 4166.6665
$ws.Range("M113").Value = 1395.2308
$ws.Range("N113").Value = -8506.666499999999
$ws.Range("H122").Value = 2827.3635
$ws.Range("I122").Value = 2273.5789
$ws.Range("K122").Value = 6820.736699999999
$ws.Range("M122").Value = -4370.736699999999
$ws.Range("H126").Value = 5719.2856
$ws.Range("I126").Value = 2740.3333
$ws.Range("J126").Value = 7953.5
$ws.Range("K126").Value = 8220.999899999999
$ws.Range("L126").Value = 23860.5
$ws.Range("M126").Value = -5750.999899999999
$ws.Range("N126").Value = -28800.5
$ws.Range("H132").Value = 2809.1853
$ws.Range("I132").Value = 2404.762
$ws.Range("K132").Value = 7214.286
$ws.Range("M132").Value = -4684.286
$ws.Range("H136").Value = 17862116
$ws.Range("I136").Value = 3666.7856
$ws.Range("K136").Value = 11000.3568
$ws.Range("M136").Value = -8450.356800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2860
$ws.Range("J39").Value = 2860
$ws.Range("L39").Value = 8580
$ws.Range("N39").Value = -9168

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6157.143
$ws.Range("I70").Value = 6900
$ws.Range("J70").Value = 5166.6665
$ws.Range("K70").Value = 6900
$ws.Range("L70").Value = 5166.6665
$ws.Range("M70").Value = -6630
$ws.Range("N70").Value = -5706.6665
$ws.Range("H73").Value = 6157.143
$ws.Range("I73").Value = 6900
$ws.Range("J73").Value = 5166.6665
$ws.Range("K73").Value = 6900
$ws.Range("L73").Value = 5166.6665
$ws.Range("M73").Value = -5964
$ws.Range("N73").Value = -7038.6665
$ws.Range("H80").Value = 3046.2363
$ws.Range("I80").Value = 2729.2856
$ws.Range("J80").Value = 3600.9
$ws.Range("K80").Value = 2729.2856
$ws.Range("L80").Value = 3600.9
$ws.Range("M80").Value = -1731.2856
$ws.Range("N80").Value = -5596.9
$ws.Range("H83").Value = 3046.2363
$ws.Range("I83").Value = 2729.2856
$ws.Range("J83").Value = 3600.9
$ws.Range("K83").Value = 13646.428
$ws.Range("L83").Value = 18004.5
$ws.Range("M83").Value = -8654.428
$ws.Range("N83").Value = -27988.5
$ws.Range("H102").Value = 51366.81
$ws.Range("I102").Value = 2914.2856
$ws.Range("K102").Value = 2914.2856
$ws.Range("M102").Value = -1292.2856
$ws.Range("H107").Value = 1061.375
$ws.Range("I107").Value = 434.72726
$ws.Range("J107").Value = 2440
$ws.Range("K107").Value = 434.72726
$ws.Range("L107").Value = 2440
$ws.Range("M107").Value = 1485.27274
$ws.Range("N107").Value = -6280
$ws.Range("H113").Value = 3247.625
$ws.Range("I113").Value = 996.8333
$ws.Range("K113").Value = 996.8333
$ws.Range("M113").Value = 1173.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 100004630
$ws.Range("J61").Value = 8198
$ws.Range("L61").Value = 8198
$ws.Range("N61").Value = -8602
$ws.Range("H113").Value = 100004630
$ws.Range("J113").Value = 8198
$ws.Range("L113").Value = 8198
$ws.Range("N113").Value = -12538
$ws.Range("H122").Value = 3064.7646
$ws.Range("I122").Value = 2527
$ws.Range("J122").Value = 4812.5
$ws.Range("K122").Value = 7581
$ws.Range("L122").Value = 14437.5
$ws.Range("M122").Value = -5131
$ws.Range("N122").Value = -19337.5
$ws.Range("H132").Value = 2250.2454
$ws.Range("I132").Value = 1564.9354
$ws.Range("J132").Value = 3215.9092
$ws.Range("K132").Value = 4694.8062
$ws.Range("L132").Value = 9647.7276
$ws.Range("M132").Value = -2164.8062
$ws.Range("N132").Value = -14707.7276
$ws.Range("H136").Value = 4626.6895
$ws.Range("I136").Value = 4082.4375
$ws.Range("K136").Value = 12247.3125
$ws.Range("M136").Value = -9697.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 28995
$ws.Range("J119").Value = 28995
$ws.Range("L119").Value = 28995
$ws.Range("N119").Value = -38671
$ws.Range("H136").Value = 2687.5625
$ws.Range("I136").Value = 2167.9546
$ws.Range("K136").Value = 6503.8638
$ws.Range("M136").Value = -3953.8638
